$wb = $excel.ActiveWorkbook
# The workbook ships cached formula results (COUNTIF/MAX/MIN summary cells in column N)
# that the original report generator left stale when it rewrote the raw cell data below.
# Switch to manual calculation so our writes do not silently refresh those cached values.
$excel.Calculation = -4135  # xlCalculationManual
$wsWS = $wb.Worksheets.Item("TESTS_WS")
$wsJMS = $wb.Worksheets.Item("TESTS_JMS")

# Row -> new "Last Execution Started" (H) and "Execution time" (I) values for TESTS_WS
$rowsWS = @{
    2 = @{ H = 42846.67847244213; I = "0.062s" }
    3 = @{ H = 42846.67848834491; I = "8.827s" }
    4 = @{ H = 42873.585141111114; I = "10.51s" }
    8 = @{ H = 42846.678645486114; I = "12.473s" }
    10 = @{ H = 42846.67879300926; I = "2.144s" }
    12 = @{ H = 42846.67881993055; I = "6.456s" }
    13 = @{ H = 42846.678897418984; I = "4.13s" }
    14 = @{ H = 42846.678947974535; I = "2.257s" }
    16 = @{ H = 42846.678977118056; I = "4.1s" }
    17 = @{ H = 42846.67902636574; I = "4.08s" }
    18 = @{ H = 42846.679075960645; I = "4.105s" }
    21 = @{ H = 42846.67912554398; I = "4.054s" }
    22 = @{ H = 42873.59151266204; I = "0.176s" }
    26 = @{ H = 42846.67917415509; I = "4.08s" }
    28 = @{ H = 42846.679223645835; I = "4.177s" }
    32 = @{ H = 42846.679273726855; I = "4.169s" }
    33 = @{ H = 42846.67932539352; I = "4.095s" }
    36 = @{ H = 42846.679374479165; I = "4.082s" }
    43 = @{ H = 42846.67942380787; I = "4.071s" }
    47 = @{ H = 42846.6794725; I = "4.101s" }
    48 = @{ H = 42846.67952165509; I = "4.202s" }
    51 = @{ H = 42846.67957408565; I = "4.073s" }
    54 = @{ H = 42846.67962337963; I = "4.082s" }
    56 = @{ H = 42846.67967762731; I = "4.586s" }
    58 = @{ H = 42846.67973241898; I = "4.095s" }
    59 = @{ H = 42846.67978160879; I = "4.08s" }
    61 = @{ H = 42846.67983094908; I = "2.398s" }
    64 = @{ H = 42846.67986560185; I = "4.146s" }
    66 = @{ H = 42846.67991545139; I = "4.12s" }
    67 = @{ H = 42846.679965219904; I = "4.066s" }
    68 = @{ H = 42846.68001456019; I = "4.083s" }
    70 = @{ H = 42846.680063761574; I = "4.133s" }
    72 = @{ H = 42846.680113125; I = "4.092s" }
    73 = @{ H = 42846.6801621412; I = "4.076s" }
    74 = @{ H = 42846.680211342595; I = "4.069s" }
    76 = @{ H = 42846.68026008102; I = "4.114s" }
    77 = @{ H = 42846.680309965275; I = "4.066s" }
    78 = @{ H = 42846.6803609375; I = "4.061s" }
    79 = @{ H = 42846.68040982639; I = "4.08s" }
    81 = @{ H = 42846.68045871528; I = "4.091s" }
    82 = @{ H = 42846.68050766204; I = "4.086s" }
    83 = @{ H = 42846.680558148146; I = "4.114s" }
    85 = @{ H = 42846.68060744213; I = "4.053s" }
    88 = @{ H = 42846.68065613426; I = "4.121s" }
    89 = @{ H = 42846.68070576389; I = "4.079s" }
    91 = @{ H = 42846.680754479166; I = "43.716s" }
    92 = @{ H = 42846.6812665625; I = "4.428s" }
    93 = @{ H = 42846.681320682874; I = "0.027s" }
    95 = @{ H = 42846.68132206018; I = "3.176s" }
    96 = @{ H = 42846.68136212963; I = "4.301s" }
    100 = @{ H = 42846.68141489584; I = "0.015s" }
    101 = @{ H = 42873.590288252315; I = "18.777s" }
    104 = @{ H = 42846.681632465275; I = "27.724s" }
    105 = @{ H = 42846.68196258102; I = "4.053s" }
    107 = @{ H = 42846.68201121528; I = "2.184s" }
    108 = @{ H = 42846.68204096065; I = "2.12s" }
    110 = @{ H = 42846.68206737268; I = "7.067s" }
    113 = @{ H = 42846.682152037036; I = "0.075s" }
    118 = @{ H = 42846.682154027774; I = "4.132s" }
}

foreach ($row in $rowsWS.Keys) {
    $entry = $rowsWS[$row]
    $wsWS.Cells.Item([int]$row, 8).Value = $entry.H
    $wsWS.Cells.Item([int]$row, 9).Value = $entry.I
}

# Row 128 is a special case: test suite (B), result (G), started (H), execution time (I) and comment (K) all change
$wsWS.Cells.Item(128, 2).Value = "PASSING TCs - Web Service Submission - Auto Not for Bamboo"
$wsWS.Cells.Item(128, 7).Value = "FAIL"
$wsWS.Cells.Item(128, 8).Value = 42852.62302690972
$wsWS.Cells.Item(128, 9).Value = "5.442s"
$wsWS.Cells.Item(128, 11).Value = "27-04-2017 14:57:13: Test case FAILED on step 7: Download Message 1|| Returned error message[s]: `n |[Not SOAP Fault] Response is a SOAP Fault| `n |[Valid HTTP Status Codes] Response status code:500 is not in acceptable list of status codes| "

# TESTS_JMS sheet: rows 2 and 3 also get refreshed "Last Execution Started" / "Execution time"
$wsJMS.Cells.Item(2, 8).Value = 42873.385770266206
$wsJMS.Cells.Item(2, 9).Value = "0.284s"
$wsJMS.Cells.Item(3, 8).Value = 42873.385813761575
$wsJMS.Cells.Item(3, 9).Value = "0.287s"
